$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 265 - this shifts the former rows 265..347
# down to 266..348 (dimension grows from R347 to R348), matching the
# weekly price-list entry that was added to the "Betarraga" sheet.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with the new weekly record.
$ws.Range("A265").Value = 4
$ws.Range("B265").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C265").Value = "Los Lagos"
$ws.Range("D265").Value = 44809
$ws.Range("E265").Value = 10
$ws.Range("F265").Value = 100114014
$ws.Range("G265").Value = "Betarraga"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 500
$ws.Range("K265").Value = 1500
$ws.Range("L265").Value = 1500
$ws.Range("M265").Value = 1500
$ws.Range("N265").Value = "`$/paquete 5 unidades"
$ws.Range("O265").Value = "Región del Maule"
$ws.Range("P265").Value = 300
$ws.Range("Q265").Value = 5
$ws.Range("R265").Value = "Hortaliza"
